$d = $word.ActiveDocument

# 1. Profile paragraph: drop the "jQuery" mention from the technology list
#    ("...Javascript, jQuery, SQL Server..." -> "...Javascript, SQL Server...").
$d.Content.Find.Execute(", jQuery", $true, $false, $false, $false, $false, `
                         $false, 1, $false, "", 2) | Out-Null

# 2. Optix Software role: rewrite the description to reflect the move to the
#    DevOps team (previously on the Till feature team) and the new focus on
#    the backend payments service.
$findText = "Software Developer working on the Till feature team on Optix 2. Optix 2 uses a microservices architecture. During my time at Optix I have been working with .Net 8, Entity Framework, Elastic Search, RabbitMq, Moq/Nsubstitute, Grafana, Docker and many other technologies."
$replaceText = "Software Developer working on the DevOps team looking at continuous improvement projects and other cross cutting concerns across the microservices estate. Before that I was part of the Till feature team on Optix 2, looking at enhancing the backend service to allow our customers to take payment from patients and sell their products and services."
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                         $false, 1, $false, $replaceText, 2) | Out-Null
